$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 4 mirrors the existing rows: barcode/id are numeric-looking
# strings that must stay text, so format as Text before entry and then
# drop back to the Normal style (matches the unstyled data rows above).
$ws.Range("A4:B4").NumberFormat = "@"
$ws.Range("A4").Value = "749309"
$ws.Range("B4").Value = "749309"
$ws.Range("A4:B4").Style = "Normal"

$ws.Range("C4").Value = "12/25/24 3:00pm"
$ws.Range("D4").Value = "test@email.edu"
$ws.Range("E4").Value = "student_class"
$ws.Range("F4").Value = "instructor"
$ws.Range("G4").Value = "name"
$ws.Range("H4").Value = "role"
$ws.Range("I4").Value = "department"
$ws.Range("J4").Value = "institution"
$ws.Range("K4").Value = "service"
$ws.Range("L4").Value = "caseName"
